$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to the player's name
$ws.Name = "Jonny Bairstow"

# Headers (row 1) - a new "matchNo" column is inserted at the front, shifting everything right
$headers = @("matchNo","teamName","batterName","states","runs","balls","fours","sixes","sr","opponentTeamName","venue","date","result")

# All values in this sheet are stored as text (even the numeric-looking ones),
# so force Text formatting before writing so Excel doesn't coerce them to numbers.
$headerRange = $ws.Range("A1:M1")
$headerRange.NumberFormat = "@"
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Full per-match batting data (rows 2-8)
$data = @(
    @("28th","Sunrisers Hyderabad","Jonny Bairstow","c Anuj Rawat b Tewatia","30","21","4","1","142.85","Rajasthan Royals","Delhi","May 02","Royals won by 55 runs"),
    @("23rd","Sunrisers Hyderabad","Jonny Bairstow","c Chahar b Curran","7","5","1","0","140.00","Chennai Super Kings","Delhi","April 28","Super Kings won by 7 wickets (with 9 balls remaining)"),
    @("20th","Sunrisers Hyderabad","Jonny Bairstow","c Dhawan b Avesh Khan","38","18","3","4","211.11","Delhi Capitals","Chennai","April 25","Match tied (Capitals won the one-over eliminator)"),
    @("3rd","Sunrisers Hyderabad","Jonny Bairstow","c Rana b Cummins","55","40","5","3","137.50","Kolkata Knight Riders","Chennai","April 11","KKR won by 10 runs"),
    @("6th","Sunrisers Hyderabad","Jonny Bairstow","c †de Villiers b Shahbaz Ahmed","12","13","1","0","92.30","Royal Challengers Bangalore","Chennai","April 14","RCB won by 6 runs"),
    @("9th","Sunrisers Hyderabad","Jonny Bairstow","hit wicket b KH Pandya","43","22","3","4","195.45","Mumbai Indians","Chennai","April 17","Mumbai won by 13 runs"),
    @("14th","Sunrisers Hyderabad","Jonny Bairstow","","63","56","3","3","112.50","Punjab Kings","Chennai","April 21","Sunrisers won by 9 wickets (with 8 balls remaining)")
)

$dataRange = $ws.Range("A2:M8")
$dataRange.NumberFormat = "@"

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}
